# Add the "expense" income rows (salary/trading/youtube/interest/salary) to
# the income sheet so it has enough data to drive a Line Chart, per the
# commit message "adding the expense section where i show all the expense
# in the form of Line Chart".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: fill in the missing Source label and update Amount/Date ---
$ws.Range("A2").Value = "salary"
$ws.Range("B2").Value = 12500
$ws.Range("C2").Value = 45809.22928240741

# --- New rows 3-6: copy the date-formatted style from C2 down first, ---
# --- so the new date cells pick up the same number format (mm-dd-yy). ---
$ws.Range("C2").Copy()
$ws.Range("C3:C6").PasteSpecial(-4122)

$ws.Range("A3").Value = "Trading"
$ws.Range("B3").Value = 10000
$ws.Range("C3").Value = 45809.22928240741

$ws.Range("A4").Value = "Youtube Revenue"
$ws.Range("B4").Value = 9500
$ws.Range("C4").Value = 45809.22928240741

$ws.Range("A5").Value = "Interest From Saving account"
$ws.Range("B5").Value = 4300
$ws.Range("C5").Value = 45806.22928240741

$ws.Range("A6").Value = "Salary"
$ws.Range("B6").Value = 5600
$ws.Range("C6").Value = 45797.92650462963
